$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatosCuenta")
$originalActive = $wb.ActiveSheet

$ws.Range("A2").Value = "SmokePre"
$ws.Range("B2").Value = "SmokePreHotFix"

$ws.Activate()
$ws.Range("D2").Select()
$originalActive.Activate()
